# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted above the existing row 31
# ("Femacal de La Calera" - Papaya), pushing the former rows 31-36 down
# to 32-37 and growing the sheet's used range from A1:T36 to A1:T37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 31; existing rows 31-36 shift to 32-37.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new weekly data point.
$ws.Range("A31").Value = 3
$ws.Range("B31").Value = "Femacal de La Calera"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44466
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100108
$ws.Range("H31").Value = "Tropicales y subtropicales"
$ws.Range("I31").Value = 100108004
$ws.Range("J31").Value = "Papaya"
$ws.Range("K31").Value = "Cultivar IV Región"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 70
$ws.Range("N31").Value = 16000
$ws.Range("O31").Value = 16000
$ws.Range("P31").Value = 16000
$ws.Range("Q31").Value = '$/bandeja 10 kilos'
$ws.Range("R31").Value = "Provincia del Elquí"
$ws.Range("S31").Value = 1600
$ws.Range("T31").Value = 10
